# Corrected excel sheets for application fix issues
#
# Applies the recorded edits to:
#   - Summary            (recomputed "Over Due"/"Outstanding" figures)
#   - Repayment schedule  (new "O" column + re-amortised schedule rows)
#   - Transactions        (renumbered transaction ids)
# and finally leaves the "Transactions" sheet as the active tab/selection,
# matching the saved workbook view state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 213.52
$wsSummary.Range("E3").Value = 113.52

# ---------------------------------------------------------------------
# 2. Repayment schedule sheet
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# -- add the new "O" column (copies the format of the adjoining "N"
#    column so the new cells carry the same style index as their
#    neighbours) --------------------------------------------------------
$wsRepay.Range("N2").Copy()
$wsRepay.Range("O2").PasteSpecial(-4122)

$wsRepay.Range("N3").Copy()
$wsRepay.Range("O3").PasteSpecial(-4122)
$wsRepay.Range("O3").Value = 0

$wsRepay.Range("N4").Copy()
$wsRepay.Range("O4").PasteSpecial(-4122)
$wsRepay.Range("O4").Value = 0

$wsRepay.Range("N5").Copy()
$wsRepay.Range("O5").PasteSpecial(-4122)
$wsRepay.Range("O5").Value = 0

$wsRepay.Range("N6").Copy()
$wsRepay.Range("O6").PasteSpecial(-4122)
$wsRepay.Range("O6").Value = 0

$wsRepay.Range("N7").Copy()
$wsRepay.Range("O7").PasteSpecial(-4122)
$wsRepay.Range("O7").Value = 0

$wsRepay.Range("N8").Copy()
$wsRepay.Range("O8").PasteSpecial(-4122)
$wsRepay.Range("O8").Value = 0

# -- re-amortised schedule values --------------------------------------
$wsRepay.Range("C4").Value = 42095
$wsRepay.Range("F4").Value = 923.19
$wsRepay.Range("G4").Value = 3212.27
$wsRepay.Range("H4").Value = 41.35

$wsRepay.Range("B5").Value = 30
$wsRepay.Range("C5").Value = 42125
$wsRepay.Range("F5").Value = 932.42
$wsRepay.Range("G5").Value = 2279.85
$wsRepay.Range("H5").Value = 32.119999999999997

$wsRepay.Range("B6").Value = 31
$wsRepay.Range("C6").Value = 42156
$wsRepay.Range("F6").Value = 941.74
$wsRepay.Range("G6").Value = 1338.11
$wsRepay.Range("H6").Value = 22.8

$wsRepay.Range("B7").Value = 30
$wsRepay.Range("C7").Value = 42186
$wsRepay.Range("F7").Value = 951.16
$wsRepay.Range("G7").Value = 386.95
$wsRepay.Range("H7").Value = 13.38

$wsRepay.Range("B8").Value = 31
$wsRepay.Range("C8").Value = 42217
$wsRepay.Range("F8").Value = 386.95
$wsRepay.Range("H8").Value = 3.87
$wsRepay.Range("K8").Value = 390.82
$wsRepay.Range("P8").Value = 390.82

# -- the sheet's own selection is no longer on L8; the saved file shows
#    a full-row selection just below the data (row 9) ------------------
$wsRepay.Range("A9:XFD9").Select()

# ---------------------------------------------------------------------
# 3. Transactions sheet
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 83
$wsTrans.Range("A3").Value = 82

# Transactions becomes the active sheet/tab, with A2:L3 selected and
# A2 as the active cell.
$wsTrans.Range("A2:L3").Select()
$wsTrans.Activate()
